$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted before the existing row 145,
# pushing the old rows 145-153 down to 146-154 and extending the used
# range to A1:T154.
$ws.Rows.Item(145).Insert()

$ws.Range("A145").Value = 9
$ws.Range("B145").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C145").Value = "Metropolitana"
$ws.Range("D145").Value = 44568
$ws.Range("E145").Value = 13
$ws.Range("F145").Value = "Fruta"
$ws.Range("G145").Value = 100101
$ws.Range("H145").Value = "Berries"
$ws.Range("I145").Value = 100101001
$ws.Range("J145").Value = "Arándano (blue)"
$ws.Range("K145").Value = "Sin especificar"
$ws.Range("L145").Value = "Primera"
$ws.Range("M145").Value = 710
$ws.Range("N145").Value = 3500
$ws.Range("O145").Value = 4000
$ws.Range("P145").Value = 3768
$ws.Range("Q145").Value = "$/bandeja 2 kilos"
$ws.Range("R145").Value = "Provincia de Linares"
$ws.Range("S145").Value = 1884
$ws.Range("T145").Value = 2
